# tabla_padron_actualizado.xlsx - "actualizamos archivos y borramos otros"
#
# For a set of rows in Sheet1, flip the "Esta En pareja? (Conyuge)" (column T)
# and/or "hijo/s" (column X) indicator cells from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToUpdate = @(
    "X2", "X4", "T5", "X5", "X6", "X7", "X10", "X11", "X13", "T15",
    "X15", "X16", "X17", "T18", "X18", "T19", "X23", "X24", "X25", "X26",
    "X27", "X28", "T32", "X32", "T33", "X33", "T34", "T35", "X35", "X36",
    "X37", "X39", "X40", "X42", "X48", "X49", "X52", "X53", "X55", "X58",
    "X61", "X62", "X64", "X65", "X66", "X67", "X68", "X70", "X74", "X76",
    "X77", "X78", "X81", "X82", "X85", "X90", "X91", "X95", "X98", "X101",
    "X102", "X103", "X104", "X105", "X107", "X108", "X109", "X112", "X114", "X118",
    "X119", "X120", "X124", "X125", "X127", "X128", "X129", "X134", "X136", "X138",
    "X158", "T194"
)

foreach ($cellRef in $cellsToUpdate) {
    $ws.Range($cellRef).Value = 1
}

# Leave the sheet's selection where the saved workbook had it.
$ws.Range("X1").Select()
